# Anonymize "fedcore" -> "approach" in the header rows, add a thin
# top/bottom (and right, for the rightmost column of each merged group)
# border under the merged title cells, and drop the stray empty G5 cell
# on the computational_comparison sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# --- Build the two border styles once (on quality_comparison!C1/D1), then
#     propagate them by copy/paste-format so every other cell that needs the
#     same look reuses the identical cell style instead of minting a fresh
#     (and instantly orphaned) one for every intermediate border edge. ---

$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.Item(8).LineStyle = 1    # xlEdgeTop
$c1.Borders.Item(9).LineStyle = 1    # xlEdgeBottom

$d1 = $ws1.Range("D1")
$c1.Copy()
$d1.PasteSpecial(-4122)              # xlPasteFormats
$d1.Borders.Item(10).LineStyle = 1   # xlEdgeRight

# --- computational_comparison: C1/D1 and F1/G1 get the same two styles ---
$c1b = $ws2.Range("C1")
$c1.Copy()
$c1b.PasteSpecial(-4122)

$d1b = $ws2.Range("D1")
$d1.Copy()
$d1b.PasteSpecial(-4122)

$f1b = $ws2.Range("F1")
$c1.Copy()
$f1b.PasteSpecial(-4122)

$g1b = $ws2.Range("G1")
$d1.Copy()
$g1b.PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Text anonymization ---
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- Drop the stray empty cell G5 on computational_comparison ---
$ws2.Range("G5").ClearContents()
